# Update evaluation workbook to latest seigot/master:
#  - branch_name column (A2:A22) "ish06b" -> "ish08a"
#  - SHAPE_LIST_MAX column (I2:I22) 1000 -> 180
#  - active cell / selection moved to I2 on the "list" sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("list")

for ($row = 2; $row -le 22; $row++) {
    $ws.Cells.Item($row, 1).Value = "ish08a"
    $ws.Cells.Item($row, 9).Value = 180
}

$ws.Range("I2").Select() | Out-Null
